$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the third data row (row 3) entirely - the sheet goes from
# A1:Z3 (header + 2 data rows) down to A1:Z2 (header + 1 data row).
$ws.Rows.Item(3).Delete()

# Update the remaining data row (row 2) with the new automation-run values.
$ws.Range("C2").Value = "AutoCustAibMr_0405427"
$ws.Range("D2").Value = "DHSTOlBIea"
$ws.Range("I2").Value = "2 yrd"

# U2 ("10548") looks numeric, so a plain .Value assignment would store it
# as a number cell. Round-trip it through a text formula + paste-as-values
# so it lands back in the sheet as a shared-string text cell (matching the
# original cell's lack of an explicit numeric style), same as the other
# text cells here.
$u2 = $ws.Range("U2")
$u2.Formula = '="10548"'
$u2.Copy()
$u2.PasteSpecial(-4163)

# Move the selection / scroll position to match the post-edit view.
$ws.Range("H6").Select()
